# Insert a new data row at row 266, shifting the existing rows 266-319 down
# to 267-320 (dimension grows from A1:T319 to A1:T320), then fill the newly
# inserted row 266 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(266).Insert()

$ws.Cells.Item(266, 1).Value = 9
$ws.Cells.Item(266, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(266, 3).Value = "Metropolitana"
$ws.Cells.Item(266, 4).Value = 44995
$ws.Cells.Item(266, 5).Value = 13
$ws.Cells.Item(266, 6).Value = "Fruta"
$ws.Cells.Item(266, 7).Value = 100101
$ws.Cells.Item(266, 8).Value = "Berries"
$ws.Cells.Item(266, 9).Value = 100101001
$ws.Cells.Item(266, 10).Value = "Arándano (blue)"
$ws.Cells.Item(266, 11).Value = "Sin especificar"
$ws.Cells.Item(266, 12).Value = "Primera"
$ws.Cells.Item(266, 13).Value = 200
$ws.Cells.Item(266, 14).Value = 3000
$ws.Cells.Item(266, 15).Value = 3000
$ws.Cells.Item(266, 16).Value = 3000
$ws.Cells.Item(266, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(266, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(266, 19).Value = 1500
$ws.Cells.Item(266, 20).Value = 2
